$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 44029
$ws.Range("B25").NumberFormat = "YYYY-MM-DD"

$ws.Range("C25").Value = 21965
$ws.Range("D25").Value = 299
$ws.Range("E25").Value = 1644
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 9.18
$ws.Range("H25").Value = 21.99

$ws.Range("J25").Value = $true

$ws.Range("K25").Value = 17913
$ws.Range("L25").Value = 291

$ws.Range("O25").Value = "Success!"
